$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws3 = $wb.Worksheets.Item("2017 LEAVE BALANCE")

# ---------------------------------------------------------------------------
# Sheet "2017 LEAVE BALANCE" (Table1) — new leave entries for rows 53-56
# (Done first so the new shared strings it introduces land before the
#  "2024" string added later on the other sheet, matching the source order.)
# ---------------------------------------------------------------------------
$ws3.Range("A53").Value = 45170
$ws3.Range("B53").Value = "SL(5-0-0)"
$ws3.Range("H53").Value = 5
$ws3.Range("K53").Value = "9/25-29/2023"

$ws3.Range("A54").Value = 45200
$ws3.Range("B54").Value = "VL(2-0-0)"
$ws3.Range("D54").Value = 2
$ws3.Range("K54").Value = "10/23,24/2023"

$ws3.Range("B55").Value = "SL(1-0-0)"
$ws3.Range("H55").Value = 1
# K55 becomes a date value (copy the date format used by K52 first).
$ws3.Range("K55").Value = 45224
$ws3.Range("K52").Copy()
$ws3.Range("K55").PasteSpecial(-4122)

$ws3.Range("A56").Value = 45261
$ws3.Range("B56").Value = "SL(2-0-0)"
$ws3.Range("H56").Value = 2
$ws3.Range("K56").Value = "12/14-15/2023"

# ---------------------------------------------------------------------------
# Sheet "2018 LEAVE CREDITS" (Table13) — extra EARNED entries + new period
# rows continuing the monthly-date pattern down to row 98, plus a "2024"
# year-header row (88) and one extra trailing table row (appended as 103).
# ---------------------------------------------------------------------------
$ws2.Range("C83").Value = 1.25

$ws2.Range("A84").Value = 45170
$ws2.Range("C84").Value = 1.25

$ws2.Range("A85").Value = 45200
$ws2.Range("C85").Value = 1.25

$ws2.Range("A86").Value = 45231
$ws2.Range("C86").Value = 1.25

$ws2.Range("A87").Value = 45261

# Row 88 is a year-header row like the existing "2018"/"2019"/... rows:
# force text, then copy the bold/quote-prefixed date format from A75 ("2023").
$ws2.Range("A88").Value = "'2024"
$ws2.Range("A75").Copy()
$ws2.Range("A88").PasteSpecial(-4122)

$ws2.Range("A89").Value = 45292
$ws2.Range("A90").Value = 45323
$ws2.Range("A91").Value = 45352
$ws2.Range("A92").Value = 45383
$ws2.Range("A93").Value = 45413
$ws2.Range("A94").Value = 45444
$ws2.Range("A95").Value = 45474
$ws2.Range("A96").Value = 45505
$ws2.Range("A97").Value = 45536
$ws2.Range("A98").Value = 45566

# ---------------------------------------------------------------------------
# Append a table row: the old last row (102, distinct border/style) moves to
# 103, and the new row 102 takes on the regular interior-row styling.
# ---------------------------------------------------------------------------
$lo2 = $ws2.ListObjects.Item(1)

# Materialise row 103 first, then pull the current (pre-edit) formatting of
# row 102 — the special "bottom of table" style — down onto it.
$ws2.Range("A103:K103").Value = 0
$ws2.Range("A102:K102").Copy()
$ws2.Range("A103:K103").PasteSpecial(-4122)
$ws2.Range("A103:K103").ClearContents()

# Now restyle row 102 itself to match the regular rows above it (row 101).
$ws2.Range("A101:K101").Copy()
$ws2.Range("A102:K102").PasteSpecial(-4122)
$ws2.Range("A102:K102").ClearContents()

# Grow the table to include the new trailing row, then restore the
# calculated-column formula on both row 102 and the new row 103.
$lo2.Resize($ws2.Range("A8:K103"))
$ws2.Range("G102").Formula = '=IF(ISBLANK(Table13[[#This Row],[EARNED]]),"",Table13[[#This Row],[EARNED]])'
$ws2.Range("G103").Formula = '=IF(ISBLANK(Table13[[#This Row],[EARNED]]),"",Table13[[#This Row],[EARNED]])'

# ---------------------------------------------------------------------------
# View state — restore each sheet's active-cell selection, leaving
# "2017 LEAVE BALANCE" as the active tab (as it was originally).
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B87").Select()

$ws3.Activate()
$ws3.Range("K56").Select()
